$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift "no_telepon" and "email" left into C1/D1 (removing "alamat" column)
$ws.Range("C1").Value = "no_telepon"
$ws.Range("D1").Value = "email"

# Clear the now-unused E1 cell (previously "email") and the old note row (C2)
$ws.Range("E1").ClearContents()
$ws.Range("C2").ClearContents()

# Update selection to match the new active cell
$ws.Range("C6").Select()

$wb.Save()
